$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E, shifting old D:K to F:M
$ws.Columns("D:E").Insert()

# Copy number formats from the (now-shifted) F column group into the new D:E columns
$ws.Range("F7:M7").Copy()
$ws.Range("D7:E7").PasteSpecial(-4122)
$ws.Range("F8:M35").Copy()
$ws.Range("D8:E35").PasteSpecial(-4122)
$ws.Range("F38:M38").Copy()
$ws.Range("D38:E38").PasteSpecial(-4122)
$ws.Range("F39:M77").Copy()
$ws.Range("D39:E77").PasteSpecial(-4122)
$ws.Range("F80:M80").Copy()
$ws.Range("D80:E80").PasteSpecial(-4122)
$ws.Range("F81:M102").Copy()
$ws.Range("D81:E102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# Populate the new column D and E values (latest two quarters) row by row
$ws.Cells.Item(7,4).Value = 43463
$ws.Cells.Item(7,5).Value = 43372
$ws.Cells.Item(8,4).Value = 384700
$ws.Cells.Item(8,5).Value = 264500
$ws.Cells.Item(9,4).Value = 198200
$ws.Cells.Item(9,5).Value = 129300
$ws.Cells.Item(10,4).Value = 186500
$ws.Cells.Item(10,5).Value = 135200
$ws.Cells.Item(12,4).Value = 37500
$ws.Cells.Item(12,5).Value = 35300
$ws.Cells.Item(13,4).Value = 0
$ws.Cells.Item(13,5).Value = 0
$ws.Cells.Item(14,4).Value = "NA"
$ws.Cells.Item(14,5).Value = "NA"
$ws.Cells.Item(15,4).Value = 300
$ws.Cells.Item(15,5).Value = 300
$ws.Cells.Item(17,4).Value = 354900
$ws.Cells.Item(17,5).Value = 227300
$ws.Cells.Item(18,4).Value = 29800
$ws.Cells.Item(18,5).Value = 37200
$ws.Cells.Item(20,4).Value = 400
$ws.Cells.Item(20,5).Value = 400
$ws.Cells.Item(21,4).Value = 39700
$ws.Cells.Item(21,5).Value = 47000
$ws.Cells.Item(22,4).Value = 0
$ws.Cells.Item(22,5).Value = 0
$ws.Cells.Item(23,4).Value = 30200
$ws.Cells.Item(23,5).Value = 37600
$ws.Cells.Item(24,4).Value = 2900
$ws.Cells.Item(24,5).Value = 5700
$ws.Cells.Item(25,4).Value = 0
$ws.Cells.Item(25,5).Value = 0
$ws.Cells.Item(26,4).Value = 27300
$ws.Cells.Item(26,5).Value = 31900
$ws.Cells.Item(27,4).Value = 27300
$ws.Cells.Item(27,5).Value = 31900
$ws.Cells.Item(28,4).Value = 0
$ws.Cells.Item(28,5).Value = 0
$ws.Cells.Item(29,4).Value = -2100
$ws.Cells.Item(29,5).Value = "NA"
$ws.Cells.Item(30,4).Value = 0
$ws.Cells.Item(30,5).Value = 0
$ws.Cells.Item(31,4).Value = 0
$ws.Cells.Item(31,5).Value = 0
$ws.Cells.Item(32,4).Value = -400
$ws.Cells.Item(32,5).Value = -400
$ws.Cells.Item(33,4).Value = 25200
$ws.Cells.Item(33,5).Value = 31900
$ws.Cells.Item(34,4).Value = 0
$ws.Cells.Item(34,5).Value = 0
$ws.Cells.Item(35,4).Value = 25200
$ws.Cells.Item(35,5).Value = 31900
$ws.Cells.Item(38,4).Value = 43463
$ws.Cells.Item(38,5).Value = 43372
$ws.Cells.Item(41,4).Value = 130400
$ws.Cells.Item(41,5).Value = 100100
$ws.Cells.Item(42,4).Value = 31600
$ws.Cells.Item(42,5).Value = 35000
$ws.Cells.Item(43,4).Value = 162200
$ws.Cells.Item(43,5).Value = 109600
$ws.Cells.Item(44,4).Value = 164600
$ws.Cells.Item(44,5).Value = 160800
$ws.Cells.Item(45,4).Value = 25700
$ws.Cells.Item(45,5).Value = 36300
$ws.Cells.Item(46,4).Value = 514400
$ws.Cells.Item(46,5).Value = 441800
$ws.Cells.Item(47,4).Value = 15100
$ws.Cells.Item(47,5).Value = 15600
$ws.Cells.Item(48,4).Value = 57000
$ws.Cells.Item(48,5).Value = 54200
$ws.Cells.Item(49,4).Value = 143200
$ws.Cells.Item(49,5).Value = 148200
$ws.Cells.Item(50,4).Value = 0
$ws.Cells.Item(50,5).Value = 0
$ws.Cells.Item(51,4).Value = 0
$ws.Cells.Item(51,5).Value = 0
$ws.Cells.Item(52,4).Value = 37200
$ws.Cells.Item(52,5).Value = 31800
$ws.Cells.Item(53,4).Value = 0
$ws.Cells.Item(53,5).Value = 0
$ws.Cells.Item(54,4).Value = 767000
$ws.Cells.Item(54,5).Value = 691600
$ws.Cells.Item(57,4).Value = 136700
$ws.Cells.Item(57,5).Value = 103100
$ws.Cells.Item(58,4).Value = 0
$ws.Cells.Item(58,5).Value = 0
$ws.Cells.Item(59,4).Value = 77000
$ws.Cells.Item(59,5).Value = 67600
$ws.Cells.Item(60,4).Value = 213800
$ws.Cells.Item(60,5).Value = 170700
$ws.Cells.Item(61,4).Value = 0
$ws.Cells.Item(61,5).Value = 0
$ws.Cells.Item(62,4).Value = 17900
$ws.Cells.Item(62,5).Value = 17300
$ws.Cells.Item(63,4).Value = 0
$ws.Cells.Item(63,5).Value = 0
$ws.Cells.Item(64,4).Value = 0
$ws.Cells.Item(64,5).Value = 0
$ws.Cells.Item(65,4).Value = 0
$ws.Cells.Item(65,5).Value = 0
$ws.Cells.Item(66,4).Value = 231600
$ws.Cells.Item(66,5).Value = 188000
$ws.Cells.Item(68,4).Value = 0
$ws.Cells.Item(68,5).Value = 0
$ws.Cells.Item(69,4).Value = 0
$ws.Cells.Item(69,5).Value = 0
$ws.Cells.Item(70,4).Value = 0
$ws.Cells.Item(70,5).Value = 0
$ws.Cells.Item(71,4).Value = 0
$ws.Cells.Item(71,5).Value = 0
$ws.Cells.Item(72,4).Value = 367000
$ws.Cells.Item(72,5).Value = 341800
$ws.Cells.Item(73,4).Value = 0
$ws.Cells.Item(73,5).Value = 0
$ws.Cells.Item(74,4).Value = 0
$ws.Cells.Item(74,5).Value = 0
$ws.Cells.Item(75,4).Value = 0
$ws.Cells.Item(75,5).Value = 0
$ws.Cells.Item(76,4).Value = 535300
$ws.Cells.Item(76,5).Value = 503600
$ws.Cells.Item(77,4).Value = 0
$ws.Cells.Item(77,5).Value = 0
$ws.Cells.Item(80,4).Value = 43463
$ws.Cells.Item(80,5).Value = 43372
$ws.Cells.Item(81,4).Value = 25200
$ws.Cells.Item(81,5).Value = 31900
$ws.Cells.Item(83,4).Value = 9500
$ws.Cells.Item(83,5).Value = 9400
$ws.Cells.Item(84,4).Value = 0
$ws.Cells.Item(84,5).Value = 0
$ws.Cells.Item(85,4).Value = 0
$ws.Cells.Item(85,5).Value = 0
$ws.Cells.Item(86,4).Value = 0
$ws.Cells.Item(86,5).Value = 0
$ws.Cells.Item(87,4).Value = 0
$ws.Cells.Item(87,5).Value = 0
$ws.Cells.Item(88,4).Value = 0
$ws.Cells.Item(88,5).Value = 0
$ws.Cells.Item(89,4).Value = 31400
$ws.Cells.Item(89,5).Value = 14000
$ws.Cells.Item(91,4).Value = -7100
$ws.Cells.Item(91,5).Value = -11000
$ws.Cells.Item(92,4).Value = 0
$ws.Cells.Item(92,5).Value = 0
$ws.Cells.Item(93,4).Value = 0
$ws.Cells.Item(93,5).Value = 0
$ws.Cells.Item(94,4).Value = -3700
$ws.Cells.Item(94,5).Value = -7700
$ws.Cells.Item(96,4).Value = 0
$ws.Cells.Item(96,5).Value = 0
$ws.Cells.Item(97,4).Value = 0
$ws.Cells.Item(97,5).Value = 0
$ws.Cells.Item(98,4).Value = 0
$ws.Cells.Item(98,5).Value = 0
$ws.Cells.Item(99,4).Value = 0
$ws.Cells.Item(99,5).Value = 0
$ws.Cells.Item(100,4).Value = 2400
$ws.Cells.Item(100,5).Value = 5300
$ws.Cells.Item(101,4).Value = 200
$ws.Cells.Item(101,5).Value = -300
$ws.Cells.Item(102,4).Value = 30300
$ws.Cells.Item(102,5).Value = 11300
